$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<five>"
$ws.Range("C2").Value = 30

# Row 3
$ws.Range("B3").Value = "<she>"
$ws.Range("C3").Value = 38

# Row 4
$ws.Range("B4").Value = "<bee>"
$ws.Range("C4").Value = 30

# Row 5
$ws.Range("B5").Value = "<find>"
$ws.Range("C5").Value = 30

# Row 6
$ws.Range("B6").Value = "<they>"
$ws.Range("C6").Value = 29

# Row 7 (column C only)
$ws.Range("C7").Value = 23

# Row 8 (column C only)
$ws.Range("C8").Value = 28

# Row 9
$ws.Range("B9").Value = "<alpha>"
$ws.Range("C9").Value = 35

# Row 10
$ws.Range("B10").Value = "<come>"
$ws.Range("C10").Value = 29

# Row 11
$ws.Range("B11").Value = "<could>"
$ws.Range("C11").Value = 30

# Row 12
$ws.Range("B12").Value = "<an>"
$ws.Range("C12").Value = 27

# Row 13 (column C only)
$ws.Range("C13").Value = 33

# Row 14
$ws.Range("B14").Value = "<have>"
$ws.Range("C14").Value = 28

# Row 15 (column C only)
$ws.Range("C15").Value = 33

# Row 16
$ws.Range("B16").Value = "<elte>"
$ws.Range("C16").Value = 31

# Row 17 (column C only)
$ws.Range("C17").Value = 31

# Row 18 (column C only)
$ws.Range("C18").Value = 34
